$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (2-6) with new interpolated/swapped values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 30.79460126660187
$ws.Range("B3").Value = 0.01318042552552552

$ws.Range("A4").Value = 37.79655091322297
$ws.Range("B4").Value = 0.02589837997997998

$ws.Range("A5").Value = 39.00197108910267
$ws.Range("B5").Value = 0.0312167972972973

$ws.Range("A6").Value = 39.66155870212616
$ws.Range("B6").Value = 0.042316103003003

# Delete the now-unused rows 7 through 34
$ws.Range("A7:B34").EntireRow.Delete()
